$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)

# 1) Give the content placeholder an explicit position/size (xfrm).
$shape.Left = 49.4999981
$shape.Top = 120.92299272598426
$shape.Width = 640.9615784031496
$shape.Height = 390.3268585937008

$tf = $shape.TextFrame
$tr = $tf.TextRange

# 2) "Download and uncompress." -> "Download and uncompress file"
#    (the trailing "." run becomes " file").
$full = $tr.Text
$dotIdx = $full.IndexOf("uncompress.") + "uncompress".Length
$dotRange = $tr.Characters($dotIdx + 1, 1)
$dotRange.Text = " file"

# 3) Add a new sub-bullet "Set source root" after "Understand project structure".
$beforeLen = $tr.Text.Length
$inserted = $tr.InsertAfter("`rSet source root")
$afterLen = $tr.Text.Length
# InsertAfter returns a range spanning the whole text body, so narrow it
# down to just the newly-added paragraph text before restyling it.
$startOfNew = $beforeLen + 2
$lengthOfNew = $afterLen - $beforeLen - 1
$newPara = $tr.Characters($startOfNew, $lengthOfNew)
$newPara.IndentLevel = 2
